# Bom_excel.xlsx — "Added price to components lists"
#
# Before: A Designator | B Footprint | C Tilausmäärä | D Value | E LCSC Part # (hyperlinks)
#         row7 is a spacer row, row3 ("100u" / pajalta?) has no hyperlink (plain "pajalta?" text)
#
# After:  A Designator | B Footprint | C Quantity | D Value | E Cost / unit (blank) | F Order link (hyperlinks)
#         spacer row7 removed, "pajalta?" -> "Pajalta" plain text in column F

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- structural edits -------------------------------------------------
# Insert a new blank column at E; old E (LCSC Part # / hyperlinks) becomes F.
$ws.Columns("E:E").Insert()

# Remove the blank spacer row (old row 7) entirely; rows below shift up.
$ws.Rows("7:7").Delete()

# --- header row ---------------------------------------------------------
$ws.Range("C1").Value = "Quantity"
$ws.Range("E1").Value = "Cost / unit"
$ws.Range("F1").Value = "Order link"

# --- "pajalta" row (formerly row3, still row3) ---------------------------
# Was a plain-text note "pajalta?" with no hyperlink; becomes "Pajalta".
$ws.Range("F3").Value = "Pajalta"

Write-Output "structural + header edits applied"
